$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2710.7046
$ws.Range("I80").Value = 1081.0588
$ws.Range("J80").Value = 3736.7778
$ws.Range("K80").Value = 3243.1764
$ws.Range("L80").Value = 11210.3334
$ws.Range("M80").Value = -2245.1764
$ws.Range("N80").Value = -13206.3334
$ws.Range("H83").Value = 2710.7046
$ws.Range("I83").Value = 1081.0588
$ws.Range("J83").Value = 3736.7778
$ws.Range("K83").Value = 9729.529200000001
$ws.Range("L83").Value = 33631.00019999999
$ws.Range("M83").Value = -4737.529200000001
$ws.Range("N83").Value = -43615.00019999999
$ws.Range("H88").Value = 1437.3182
$ws.Range("J88").Value = 502.5
$ws.Range("L88").Value = 502.5
$ws.Range("N88").Value = -1314.5
$ws.Range("H91").Value = 1437.3182
$ws.Range("J91").Value = 502.5
$ws.Range("L91").Value = 502.5
$ws.Range("N91").Value = -3310.5
$ws.Range("H112").Value = 2018.619
$ws.Range("J112").Value = 2104.7896
$ws.Range("L112").Value = 6314.3688
$ws.Range("N112").Value = -8530.3688
$ws.Range("H116").Value = 4579.1113
$ws.Range("I116").Value = 4101.857
$ws.Range("K116").Value = 4101.857
$ws.Range("M116").Value = -659.857
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13368.274
$ws.Range("I32").Value = 10884.382
$ws.Range("J32").Value = 32884.57
$ws.Range("K32").Value = 10884.382
$ws.Range("L32").Value = 32884.57
$ws.Range("M32").Value = -10597.382
$ws.Range("N32").Value = -33458.57
$ws.Range("H61").Value = 3912.2632
$ws.Range("I61").Value = 3155.7856
$ws.Range("J61").Value = 6030.4
$ws.Range("K61").Value = 3155.7856
$ws.Range("L61").Value = 6030.4
$ws.Range("M61").Value = -2943.7856
$ws.Range("N61").Value = -6454.4
$ws.Range("H74").Value = 2956.9583
$ws.Range("I74").Value = 2927.0952
$ws.Range("J74").Value = 3166
$ws.Range("K74").Value = 2927.0952
$ws.Range("L74").Value = 3166
$ws.Range("M74").Value = -2053.0952
$ws.Range("N74").Value = -4914
$ws.Range("H77").Value = 2956.9583
$ws.Range("I77").Value = 2927.0952
$ws.Range("J77").Value = 3166
$ws.Range("K77").Value = 14635.476
$ws.Range("L77").Value = 15830
$ws.Range("M77").Value = -10267.476
$ws.Range("N77").Value = -24566
$ws.Range("H88").Value = 1682.0625
$ws.Range("I88").Value = 1775.1177
$ws.Range("J88").Value = 1576.6
$ws.Range("K88").Value = 1775.1177
$ws.Range("L88").Value = 1576.6
$ws.Range("M88").Value = -1369.1177
$ws.Range("N88").Value = -2388.6
$ws.Range("H91").Value = 1682.0625
$ws.Range("I91").Value = 1775.1177
$ws.Range("J91").Value = 1576.6
$ws.Range("K91").Value = 1775.1177
$ws.Range("L91").Value = 1576.6
$ws.Range("M91").Value = -371.1177
$ws.Range("N91").Value = -4384.6
$ws.Range("H132").Value = 3245.6309
$ws.Range("I132").Value = 3049.2407
$ws.Range("J132").Value = 4209.727
$ws.Range("K132").Value = 9147.722099999999
$ws.Range("L132").Value = 12629.181
$ws.Range("M132").Value = -6617.722099999999
$ws.Range("N132").Value = -17689.181
$ws.Range("H136").Value = 3912.2632
$ws.Range("I136").Value = 3155.7856
$ws.Range("J136").Value = 6030.4
$ws.Range("K136").Value = 9467.356800000001
$ws.Range("L136").Value = 18091.2
$ws.Range("M136").Value = -6917.356800000001
$ws.Range("N136").Value = -23191.2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3373.4878
$ws.Range("I134").Value = 2340.8147
$ws.Range("J134").Value = 5365.0713
$ws.Range("K134").Value = 7022.4441
$ws.Range("L134").Value = 16095.2139
$ws.Range("M134").Value = -4487.4441
$ws.Range("N134").Value = -21165.2139
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 48540
$ws.Range("I86").Value = 63178.6
$ws.Range("J86").Value = 30241.75
$ws.Range("K86").Value = 63178.6
$ws.Range("L86").Value = 30241.75
$ws.Range("M86").Value = -62055.6
$ws.Range("N86").Value = -32487.75
$ws.Range("H89").Value = 48540
$ws.Range("I89").Value = 63178.6
$ws.Range("J89").Value = 30241.75
$ws.Range("K89").Value = 315893
$ws.Range("L89").Value = 151208.75
$ws.Range("M89").Value = -310277
$ws.Range("N89").Value = -162440.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 202.4
$ws.Range("I33").Value = 102.75
$ws.Range("J33").Value = 601
$ws.Range("K33").Value = 616.5
$ws.Range("L33").Value = 3606
$ws.Range("M33").Value = -333.5
$ws.Range("N33").Value = -4172
$ws.Range("H70").Value = 227712
$ws.Range("I70").Value = 5975
$ws.Range("K70").Value = 17925
$ws.Range("M70").Value = -17610
$ws.Range("H73").Value = 227712
$ws.Range("I73").Value = 5975
$ws.Range("K73").Value = 17925
$ws.Range("M73").Value = -16833
$ws.Range("H121").Value = 2302.2964
$ws.Range("I121").Value = 253.2
$ws.Range("J121").Value = 2768
$ws.Range("K121").Value = 759.5999999999999
$ws.Range("L121").Value = 8304
$ws.Range("M121").Value = 550.4000000000001
$ws.Range("N121").Value = -10924
$ws.Range("H131").Value = 3992.087
$ws.Range("I131").Value = 4386.5
$ws.Range("J131").Value = 3852.8823
$ws.Range("K131").Value = 13159.5
$ws.Range("L131").Value = 11558.6469
$ws.Range("M131").Value = -8119.5
$ws.Range("N131").Value = -21638.6469
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 634.1429000000001
$ws.Range("I9").Value = 415.66666
$ws.Range("J9").Value = 798
$ws.Range("K9").Value = 415.66666
$ws.Range("L9").Value = 798
$ws.Range("M9").Value = -245.66666
$ws.Range("N9").Value = -1138
$ws.Range("H122").Value = 2935.15
$ws.Range("J122").Value = 4397.8
$ws.Range("L122").Value = 13193.4
$ws.Range("N122").Value = -18093.4
$ws.Range("H132").Value = 5018.96
$ws.Range("I132").Value = 4526.3335
$ws.Range("K132").Value = 13579.0005
$ws.Range("M132").Value = -11049.0005
$ws.Range("H136").Value = 15461.944
$ws.Range("J136").Value = 15461.944
$ws.Range("L136").Value = 46385.83199999999
$ws.Range("N136").Value = -51485.83199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8593.518
$ws.Range("I40").Value = 9993.071
$ws.Range("K40").Value = 9993.071
$ws.Range("M40").Value = -9857.071
$ws.Range("H122").Value = 4334.6177
$ws.Range("I122").Value = 3481.8965
$ws.Range("K122").Value = 10445.6895
$ws.Range("M122").Value = -7995.6895
$ws.Range("H136").Value = 4826.4863
$ws.Range("I136").Value = 2947.4666
$ws.Range("J136").Value = 6107.636
$ws.Range("K136").Value = 8842.399800000001
$ws.Range("L136").Value = 18322.908
$ws.Range("M136").Value = -6292.399800000001
$ws.Range("N136").Value = -23422.908
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1359.0769
$ws.Range("I132").Value = 1029.826
$ws.Range("K132").Value = 3089.478
$ws.Range("M132").Value = -559.4780000000001
